$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F4").Value = 99
$ws.Range("F5").Value = 98
$ws.Range("F6").Value = 99
$ws.Range("F7").Value = 99
$ws.Range("F8").Value = 98

$ws.Range("F9").Select()
